# Edit the "music.xlsx" workbook:
#  - Replace the placeholder "?" values with "tbd" in column E (rows 5,6,7,8,12,14)
#  - Move the active cell selection from J15 to E15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Paused?" column cells that were previously "?" to "tbd"
$cells = @("E5", "E6", "E7", "E8", "E12", "E14")
foreach ($cell in $cells) {
    $ws.Range($cell).Value = "tbd"
}

# Update the active selection to E15
$ws.Range("E15").Select()
